# "fix mistake in image 3"
#
# Slide 3 contains a sequence of repeated "Byte layout" diagrams
# (arrow + table + caption). The 3rd diagram on the slide ("image 3")
# had its table header mistakenly left as "Byte 4" instead of "New".
# While fixing that, also clean up two nearby "Message.byte(n) = "
# captions whose text had been split across two runs (one of them
# flagged with a stale spell-check err="1") on slide 2 and slide 3 -
# retyping them merges the runs and drops the stale flag.

$p = $ppt.ActivePresentation

function Fix-MessageByteCaption($slide) {
    $shape = $slide.Shapes.Item("Textfeld 33")
    $tr = $shape.TextFrame.TextRange
    # Text layout is: "Message.byte" (run1, err="1") + "(n) = " (run2) + "New" (run3)
    # Characters(1,12) == "Message.byte", Characters(13,6) == "(n) = "
    $prefix = $tr.Characters(1, 12).Text
    $rest = $tr.Characters(13, 6).Text
    if ($prefix -eq "Message.byte") {
        # Rewrite run2's text to the full merged caption (keeps run2's
        # clean rPr, without the err="1" proofing flag), then clear out
        # the now-redundant leading run.
        $tr.Characters(13, 6).Text = $prefix + $rest
        $tr.Characters(1, 12).Text = ""
    }
}

# Slide 2: "Message.byte(6) = New"
Fix-MessageByteCaption($p.Slides.Item(2))

# Slide 3: "Message.byte(7) = New" and the "Byte 4" -> "New" table typo.
$slide3 = $p.Slides.Item(3)
Fix-MessageByteCaption($slide3)

$table3 = $slide3.Shapes.Item("Tabelle 25")
$cell = $table3.Table.Cell(1, 1)
if ($cell.Shape.TextFrame.TextRange.Text -eq "Byte 4") {
    $cell.Shape.TextFrame.TextRange.Text = "New"
}
